# Add a new worksheet named "Loading" as the last sheet in the workbook.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Loading"

# Populate column B with the CSS for a loading overlay / spinner.
$ws.Range("B2").Value = '/* 遮罩层样式 */'
$ws.Range("B16").Value = '/* 转圈动画样式 */'
$ws.Range("B3").Value = '.loadingOverlay {'
$ws.Range("B4").Value = '  position: fixed;'
$ws.Range("B5").Value = '  top: 0;'
$ws.Range("B6").Value = '  left: 0;'
$ws.Range("B7").Value = '  width: 100vw;'
$ws.Range("B8").Value = '  height: 100vh;'
$ws.Range("B9").Value = '  background-color: rgba(200, 200, 200, 0.6); /* 淡灰色半透明 */'
$ws.Range("B10").Value = '  display: flex;'
$ws.Range("B11").Value = '  justify-content: center;'
$ws.Range("B12").Value = '  align-items: center;'
$ws.Range("B13").Value = '  z-index: 9999; /* 保证在最上层 */'
$ws.Range("B17").Value = '.spinner {'
$ws.Range("B18").Value = '  width: 60px;'
$ws.Range("B19").Value = '  height: 60px;'
$ws.Range("B20").Value = '  border: 6px solid #ccc;'
$ws.Range("B21").Value = '  border-top-color: #4CAF50;'
$ws.Range("B22").Value = '  border-radius: 50%;'
$ws.Range("B23").Value = '  animation: spin 1s linear infinite;'
$ws.Range("B26").Value = '@keyframes spin {'
$ws.Range("B27").Value = '  to { transform: rotate(360deg); }'
$ws.Range("B14").Value = '}'
$ws.Range("B24").Value = '}'
$ws.Range("B28").Value = '}'

# Match the selection recorded in the workbook for the new sheet.
$ws.Range("B2:B28").Select()
